# testing and addition of buffer of 1 month before new pregnancy in contraception
#
# For a set of parameter rows on the "parameter_values" sheet, insert a new
# first data point (0.2) into column B, pushing whatever was already in B
# (and, for row 33, also C) one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33 is special: it already has two values (B33, C33), so both need to
# shift right (C33 -> D33, B33 -> C33) before the new 0.2 is written to B33.
# The destination D33 cell was previously an empty, styled placeholder
# (s="2"); once it receives a real value the placeholder formatting is
# cleared, matching how Excel drops that formatting when the cell becomes a
# genuine data cell.
$ws.Range("D33").Value = $ws.Range("C33").Value2
$ws.Range("D33").ClearFormats()
$ws.Range("C33").Value = $ws.Range("B33").Value2
$ws.Range("B33").Value = 0.2

# Remaining rows: simple single-value shift, B -> C, then B = 0.2.
$rows = @(11, 25, 27, 32, 39, 40, 41, 42, 43, 55, 56, 59, 66, 67, 68)
foreach ($r in $rows) {
    $ws.Range("C$r").Value = $ws.Range("B$r").Value2
    $ws.Range("B$r").Value = 0.2
}

# Update the saved view state (scroll position + active selection) to match.
$win = $excel.ActiveWindow
$win.ScrollRow = 52
$win.ScrollColumn = 1
$ws.Range("F59").Select() | Out-Null
